$p = $ppt.ActivePresentation

# --- Slide 1: "Subtitle 2" placeholder ---------------------------------
# 1) Turn on "Shrink text on overflow" (normAutofit) now that a 4th line
#    of text is being added to the subtitle.
# 2) Append a new paragraph with the presenter's email address after the
#    "sorbonne center for artificial intelligence (SCAI)" line.
$slide1 = $p.Slides.Item(1)
$subtitle = $slide1.Shapes.Item(2)
$subtitle.TextFrame.AutoSize = 2
$subtitle.TextFrame.TextRange.InsertAfter("`rjames.gawley@gmail.com") | Out-Null

# --- Slide 2: "Content Placeholder 2" -----------------------------------
# Expand the "Install python" bullet to mention jupyter notebooks and git.
# (The intermediate assignment avoids the host's longest-common-prefix
# run-splitting so the paragraph ends up as a single run, matching how
# PowerPoint collapses a fully-retyped line back into one run.)
$slide2 = $p.Slides.Item(2)
$content = $slide2.Shapes.Item(2)
$installPara = $content.TextFrame.TextRange.Paragraphs(1, 1)
$installPara.Text = " "
$installPara.Text = "Install python, jupyter notebooks, git"
